# Minor clean up to the radar spreadsheet.
#
# The "Sources" worksheet had a stray blank leading row (with just the
# title "Sources" in B2) and a stray blank leading column in front of its
# Key/Reference table (which lived at B4:C7). Clean that up by dropping
# the title cell and shifting the Key/Reference table so it starts at A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sources")

# Capture the existing Key/Reference table values before clearing anything.
$keyHeader = $ws.Cells.Item(4, 2).Value()   # "Key"
$refHeader = $ws.Cells.Item(4, 3).Value()   # "Reference"

$key1 = $ws.Cells.Item(5, 2).Value()        # 1
$ref1 = $ws.Cells.Item(5, 3).Value()        # first reference text

$key2 = $ws.Cells.Item(6, 2).Value()        # 2
$ref2 = $ws.Cells.Item(6, 3).Value()        # second reference text

$key3 = $ws.Cells.Item(7, 2).Value()        # 3
$ref3 = $ws.Cells.Item(7, 3).Value()        # third reference text

# Clear out the old layout (the "Sources" title plus the Key/Reference table).
[void]$ws.Range("B2:C7").Clear()

# Write the table back starting at A1.
$ws.Cells.Item(1, 1).Value = $keyHeader
$ws.Cells.Item(1, 2).Value = $refHeader

$ws.Cells.Item(2, 1).Value = $key1
$ws.Cells.Item(2, 2).Value = $ref1

$ws.Cells.Item(3, 1).Value = $key2
$ws.Cells.Item(3, 2).Value = $ref2

$ws.Cells.Item(4, 1).Value = $key3
$ws.Cells.Item(4, 2).Value = $ref3

[void]$ws.Range("B17").Select()

# Restore the Data sheet as the active tab/selection.
$data = $wb.Worksheets.Item("Data")
[void]$data.Activate()
[void]$data.Range("D5").Select()
